$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Set numeric / text values (row-content permutation per commit diff) ---
$ws.Range("A6").Value = 130873722
$ws.Range("Q6").Value = 438858
$ws.Range("R6").Value = 6795118
$ws.Range("A7").Value = 130873728
$ws.Range("B7").Value = 79244
$ws.Range("D7").Value = "NT"
$ws.Range("E7").Value = 6425
$ws.Range("F7").Value = "Garnlav"
$ws.Range("G7").Value = "Alectoria sarmentosa"
$ws.Range("H7").Value = "(Ach.) Ach."
$ws.Range("Q7").Value = 438641
$ws.Range("R7").Value = 6795153
$ws.Range("A8").Value = 130873733
$ws.Range("Q8").Value = 438651
$ws.Range("R8").Value = 6795214
$ws.Range("A9").Value = 130873730
$ws.Range("Q9").Value = 438606
$ws.Range("R9").Value = 6795190
$ws.Range("A10").Value = 130873741
$ws.Range("Q10").Value = 438767
$ws.Range("R10").Value = 6795135
$ws.Range("A11").Value = 130873693
$ws.Range("B11").Value = 57073
$ws.Range("D11").Value = "LC"
$ws.Range("E11").Value = 100138
$ws.Range("F11").Value = "Tjäder"
$ws.Range("G11").Value = "Tetrao urogallus"
$ws.Range("H11").Value = "Linnaeus, 1758"
$ws.Range("M11").Value = "färsk spillning"
$ws.Range("Q11").Value = 438755
$ws.Range("R11").Value = 6795183
$ws.Range("A33").Value = 130873703
$ws.Range("B33").Value = 8451
$ws.Range("D33").Value = "LC"
$ws.Range("E33").Value = 106545
$ws.Range("F33").Value = "Mindre märgborre"
$ws.Range("G33").Value = "Tomicus minor"
$ws.Range("H33").Value = "(Hartig, 1834)"
$ws.Range("M33").Value = "äldre gnagspår"
$ws.Range("Q33").Value = 439003
$ws.Range("R33").Value = 6795150
$ws.Range("AX33").Value = "Eva Löfqvist"
$ws.Range("A34").Value = 130873697
$ws.Range("B34").Value = 91830
$ws.Range("E34").Value = 5442
$ws.Range("F34").Value = "Tallticka"
$ws.Range("G34").Value = "Porodaedalea pini"
$ws.Range("H34").Value = "(Brot.) Murrill"
$ws.Range("Q34").Value = 438905
$ws.Range("R34").Value = 6795075
$ws.Range("AX34").Value = "Eva Löfqvist, Alfhild Sehlin"
$ws.Range("A35").Value = 130873742
$ws.Range("Q35").Value = 438980
$ws.Range("R35").Value = 6795131
$ws.Range("A36").Value = 130873727
$ws.Range("B36").Value = 79244
$ws.Range("D36").Value = "NT"
$ws.Range("E36").Value = 6425
$ws.Range("F36").Value = "Garnlav"
$ws.Range("G36").Value = "Alectoria sarmentosa"
$ws.Range("H36").Value = "(Ach.) Ach."
$ws.Range("Q36").Value = 438675
$ws.Range("R36").Value = 6795125

# --- Clear cells that should no longer be present ---
$ws.Range("J6").Value = ""
$ws.Range("K6").Value = ""
$ws.Range("N6").Value = ""
$ws.Range("AF6").Value = ""
$ws.Range("K7").Value = ""
$ws.Range("L7").Value = ""
$ws.Range("M7").Value = ""
$ws.Range("N7").Value = ""
$ws.Range("J36").Value = ""
$ws.Range("K36").Value = ""
$ws.Range("L36").Value = ""
$ws.Range("M36").Value = ""
$ws.Range("N36").Value = ""
$ws.Range("AF36").Value = ""

# --- Re-create blank placeholder cells (present, but empty) ---
$ws.Range("J9").NumberFormat = "General"
$ws.Range("K9").NumberFormat = "General"
$ws.Range("N9").NumberFormat = "General"
$ws.Range("AF9").NumberFormat = "General"
$ws.Range("K11").NumberFormat = "General"
$ws.Range("L11").NumberFormat = "General"
$ws.Range("N11").NumberFormat = "General"
$ws.Range("L33").NumberFormat = "General"
$ws.Range("J34").NumberFormat = "General"
$ws.Range("K34").NumberFormat = "General"
$ws.Range("N34").NumberFormat = "General"
$ws.Range("AF34").NumberFormat = "General"
